$wb = $excel.ActiveWorkbook

# The "想去人数" (want-to-go count) column F was updated for both the
# "展览" and "全部类型" sheets, which contain mirrored data.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 365
    $ws.Range("F4").Value = 1563
    $ws.Range("F5").Value = 6
    $ws.Range("F7").Value = 395
    $ws.Range("F10").Value = 416
}
